$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-13 03:00:47"
$wsZh.Range("H2").Value = "2016-03-13 03:01:07"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-13 03:00:51"
$wsDe.Range("H2").Value = "2016-03-13 03:01:13"
